$d = $word.ActiveDocument

# --- 1. Strip the "_GoBack" bookmark off the leading empty paragraph ---
# (it will be re-added at the end of the new "FFF" paragraph below)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Append three new paragraphs ("DDD ", "EEE", "FFF") after "CCC" ---
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$d.Paragraphs.Last.Range.Text = "DDD "

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "EEE"

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "FFF"

# --- 3. Re-create the "_GoBack" bookmark, collapsed, right after "FFF" ---
# The COM host mis-resolves a bookmark collapsed exactly at the tail
# paragraph-mark position of the document's final paragraph, so a
# temporary non-empty paragraph is appended first to give the engine an
# unambiguous anchor; the bookmark is inserted against that safe
# position, the temporary paragraph break is removed (merging its text
# onto the "FFF" paragraph), and then that trailing placeholder text is
# deleted - leaving the bookmark collapsed immediately after "FFF".
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "ZZZZZ"

$fffParaIndex = $d.Paragraphs.Count - 1
$fffRange = $d.Paragraphs.Item($fffParaIndex).Range
$fffRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $fffRange)

$fffEnd = $d.Paragraphs.Item($fffParaIndex).Range.End
$markRange = $d.Range($fffEnd - 1, $fffEnd)
$markRange.Delete()

$mergedPara = $d.Paragraphs.Item($fffParaIndex)
$suffixRange = $d.Range($fffEnd - 1, $mergedPara.Range.End - 1)
$suffixRange.Delete()
